# Updates numeric price/profit figures (columns H-N) across the job sheets
# to match the refreshed market-board snapshot pulled by the scheduled runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 98.2
$ws.Range("I2").Value = 96
$ws.Range("K2").Value = 96
$ws.Range("M2").Value = 17
$ws.Range("H69").Value = 10013.808
$ws.Range("J69").Value = 10015
$ws.Range("L69").Value = 30045
$ws.Range("N69").Value = -31793
$ws.Range("H72").Value = 10013.808
$ws.Range("J72").Value = 10015
$ws.Range("L72").Value = 90135
$ws.Range("N72").Value = -98871
$ws.Range("H76").Value = 5155
$ws.Range("I76").Value = 5193.75
$ws.Range("J76").Value = 5000
$ws.Range("K76").Value = 5193.75
$ws.Range("L76").Value = 5000
$ws.Range("M76").Value = -4878.75
$ws.Range("N76").Value = -5630
$ws.Range("H79").Value = 5155
$ws.Range("I79").Value = 5193.75
$ws.Range("J79").Value = 5000
$ws.Range("K79").Value = 5193.75
$ws.Range("L79").Value = 5000
$ws.Range("M79").Value = -4101.75
$ws.Range("N79").Value = -7184
$ws.Range("H92").Value = 31251146
$ws.Range("I92").Value = 35715332
$ws.Range("J92").Value = 1839.5
$ws.Range("K92").Value = 35715332
$ws.Range("L92").Value = 1839.5
$ws.Range("M92").Value = -35714084
$ws.Range("N92").Value = -4335.5
$ws.Range("H103").Value = 807.05554
$ws.Range("I103").Value = 460.21054
$ws.Range("J103").Value = 1194.7059
$ws.Range("K103").Value = 1380.63162
$ws.Range("L103").Value = 3584.1177
$ws.Range("M103").Value = -794.6316199999999
$ws.Range("N103").Value = -4756.1177
$ws.Range("H113").Value = 22414
$ws.Range("I113").Value = 24483
$ws.Range("J113").Value = 10000
$ws.Range("K113").Value = 24483
$ws.Range("L113").Value = 10000
$ws.Range("M113").Value = -21229
$ws.Range("N113").Value = -16508
$ws.Range("H116").Value = 16698.857
$ws.Range("I116").Value = 36719
$ws.Range("K116").Value = 36719
$ws.Range("M116").Value = -33277
$ws.Range("H132").Value = 1706.7838
$ws.Range("J132").Value = 2269.75
$ws.Range("L132").Value = 6809.25
$ws.Range("N132").Value = -11869.25
$ws.Range("H137").Value = 57483.934
$ws.Range("I137").Value = 71563.25
$ws.Range("K137").Value = 214689.75
$ws.Range("M137").Value = -212139.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2937.1875
$ws.Range("I2").Value = 2149.75
$ws.Range("J2").Value = 3724.625
$ws.Range("K2").Value = 2149.75
$ws.Range("L2").Value = 3724.625
$ws.Range("M2").Value = -2036.75
$ws.Range("N2").Value = -3950.625
$ws.Range("H24").Value = 86168.8
$ws.Range("J24").Value = 86168.8
$ws.Range("L24").Value = 86168.8
$ws.Range("N24").Value = -86916.8
$ws.Range("H32").Value = 26963074
$ws.Range("I32").Value = 30569288
$ws.Range("K32").Value = 30569288
$ws.Range("M32").Value = -30569001
$ws.Range("H61").Value = 4005.92
$ws.Range("I61").Value = 3472.7273
$ws.Range("J61").Value = 4424.857
$ws.Range("K61").Value = 3472.7273
$ws.Range("L61").Value = 4424.857
$ws.Range("M61").Value = -3260.7273
$ws.Range("N61").Value = -4848.857
$ws.Range("H82").Value = 79996.664
$ws.Range("J82").Value = 79996.664
$ws.Range("L82").Value = 79996.664
$ws.Range("N82").Value = -80718.664
$ws.Range("H85").Value = 79996.664
$ws.Range("J85").Value = 79996.664
$ws.Range("L85").Value = 79996.664
$ws.Range("N85").Value = -82492.664
$ws.Range("H86").Value = 79999
$ws.Range("J86").Value = 79999
$ws.Range("L86").Value = 79999
$ws.Range("N86").Value = -82371
$ws.Range("H89").Value = 79999
$ws.Range("J89").Value = 79999
$ws.Range("L89").Value = 239997
$ws.Range("N89").Value = -251853
$ws.Range("H92").Value = 56160
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 56160
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 56160
$ws.Range("M92").ClearContents()
$ws.Range("N92").Value = -61152
$ws.Range("H100").Value = 86168.8
$ws.Range("J100").Value = 86168.8
$ws.Range("L100").Value = 86168.8
$ws.Range("N100").Value = -88332.8
$ws.Range("H116").Value = 2937.1875
$ws.Range("I116").Value = 2149.75
$ws.Range("J116").Value = 3724.625
$ws.Range("K116").Value = 2149.75
$ws.Range("L116").Value = 3724.625
$ws.Range("M116").Value = 144.25
$ws.Range("N116").Value = -8312.625
$ws.Range("H132").Value = 3297.2812
$ws.Range("I132").Value = 3146.75
$ws.Range("K132").Value = 9440.25
$ws.Range("M132").Value = -6910.25
$ws.Range("H136").Value = 4005.92
$ws.Range("I136").Value = 3472.7273
$ws.Range("J136").Value = 4424.857
$ws.Range("K136").Value = 10418.1819
$ws.Range("L136").Value = 13274.571
$ws.Range("M136").Value = -7868.1819
$ws.Range("N136").Value = -18374.571

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2937.1875
$ws.Range("I3").Value = 2149.75
$ws.Range("J3").Value = 3724.625
$ws.Range("K3").Value = 2149.75
$ws.Range("L3").Value = 3724.625
$ws.Range("M3").Value = -2035.75
$ws.Range("N3").Value = -3952.625
$ws.Range("H134").Value = 2566346.5
$ws.Range("I134").Value = 2900374.2
$ws.Range("K134").Value = 8701122.600000001
$ws.Range("M134").Value = -8698587.600000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 3926.4517
$ws.Range("I132").Value = 3515.3333
$ws.Range("K132").Value = 10545.9999
$ws.Range("M132").Value = -8015.999899999999
$ws.Range("H134").Value = 3636.6667
$ws.Range("J134").Value = 3949
$ws.Range("L134").Value = 11847
$ws.Range("N134").Value = -16917

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 1000111.4
$ws.Range("J12").Value = 1111224
$ws.Range("L12").Value = 3333672
$ws.Range("N12").Value = -3334018
$ws.Range("H139").Value = 6994.5454
$ws.Range("I139").Value = 1735
$ws.Range("K139").Value = 5205
$ws.Range("M139").Value = -65

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2844.077
$ws.Range("I102").Value = 2441.4443
$ws.Range("J102").Value = 3750
$ws.Range("K102").Value = 2441.4443
$ws.Range("L102").Value = 3750
$ws.Range("M102").Value = -819.4443000000001
$ws.Range("N102").Value = -6994
$ws.Range("H113").Value = 44374.668
$ws.Range("J113").Value = 62699.2
$ws.Range("L113").Value = 62699.2
$ws.Range("N113").Value = -67039.2
$ws.Range("H122").Value = 1569.6666
$ws.Range("I122").Value = 1522.6666
$ws.Range("J122").Value = 1616.6666
$ws.Range("K122").Value = 4567.9998
$ws.Range("L122").Value = 4849.9998
$ws.Range("M122").Value = -2117.9998
$ws.Range("N122").Value = -9749.9998
$ws.Range("H126").Value = 3724.7778
$ws.Range("J126").Value = 3953.5
$ws.Range("L126").Value = 11860.5
$ws.Range("N126").Value = -16800.5
$ws.Range("H132").Value = 3654
$ws.Range("I132").Value = 3685.1428
$ws.Range("J132").Value = 3000
$ws.Range("K132").Value = 11055.4284
$ws.Range("L132").Value = 9000
$ws.Range("M132").Value = -8525.428400000001
$ws.Range("N132").Value = -14060

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 10698
$ws.Range("I7").Value = 11300
$ws.Range("K7").Value = 11300
$ws.Range("M7").Value = -11188
$ws.Range("H16").Value = 2789.182
$ws.Range("I16").Value = 1763.3334
$ws.Range("J16").Value = 4020.2
$ws.Range("K16").Value = 1763.3334
$ws.Range("L16").Value = 4020.2
$ws.Range("M16").Value = -1593.3334
$ws.Range("N16").Value = -4360.2
$ws.Range("H40").Value = 3584.1428
$ws.Range("I40").Value = 2723.5
$ws.Range("J40").Value = 4731.6665
$ws.Range("K40").Value = 2723.5
$ws.Range("L40").Value = 4731.6665
$ws.Range("M40").Value = -2587.5
$ws.Range("N40").Value = -5003.6665
$ws.Range("H82").Value = 1485.2
$ws.Range("I82").Value = 1532.875
$ws.Range("J82").Value = 1430.7142
$ws.Range("K82").Value = 1532.875
$ws.Range("L82").Value = 1430.7142
$ws.Range("M82").Value = -1171.875
$ws.Range("N82").Value = -2152.7142
$ws.Range("H85").Value = 1485.2
$ws.Range("I85").Value = 1532.875
$ws.Range("J85").Value = 1430.7142
$ws.Range("K85").Value = 1532.875
$ws.Range("L85").Value = 1430.7142
$ws.Range("M85").Value = -284.875
$ws.Range("N85").Value = -3926.7142
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 2000
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = 2000
$ws.Range("M100").ClearContents()
$ws.Range("N100").Value = -3082
$ws.Range("H125").Value = 93994.5
$ws.Range("J125").Value = 93994.5
$ws.Range("L125").Value = 93994.5
$ws.Range("N125").Value = -103834.5
$ws.Range("H126").Value = 10698
$ws.Range("I126").Value = 11300
$ws.Range("K126").Value = 33900
$ws.Range("M126").Value = -31430
$ws.Range("H136").Value = 8483.308000000001
$ws.Range("I136").Value = 4875.75
$ws.Range("K136").Value = 14627.25
$ws.Range("M136").Value = -12077.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 11361.25
$ws.Range("I122").Value = 11814.333
$ws.Range("J122").Value = 10002
$ws.Range("K122").Value = 35442.999
$ws.Range("L122").Value = 30006
$ws.Range("M122").Value = -32992.999
$ws.Range("N122").Value = -34906
$ws.Range("H132").Value = 1970.4286
$ws.Range("I132").Value = 1798.8334
$ws.Range("J132").Value = 3000
$ws.Range("K132").Value = 5396.5002
$ws.Range("L132").Value = 9000
$ws.Range("M132").Value = -2866.5002
$ws.Range("N132").Value = -14060
$ws.Range("H136").Value = 12501956
$ws.Range("I136").Value = 1954.3636
$ws.Range("K136").Value = 5863.0908
$ws.Range("M136").Value = -3313.0908
